# GearSwap v0.812 - More bugfixing
# Add world.real_weather / world.real_weather_element variables (weather that
# ignores scholar weather-changing spells/buffs), and clarify that the
# existing world.weather / world.weather_element respect those spells/buffs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank rows right after the "weatherelement" row (old row 236)
# and before the "moon" row (old row 237). Because these rows sit inside the
# merged "World" header cell (A230:A241), Excel auto-grows that merge to
# A230:A243 as part of the insert.
$ws.Range("A237:A238").EntireRow.Insert()

# Clarify the existing weather-name / weather-element descriptions.
$ws.Cells.Item(235, 6).Value = "Current weather's name. Respects scholar spells/buffs."
$ws.Cells.Item(236, 6).Value = "Current weather's element. Respects scholar spells/buffs."

# Fill in the two newly inserted rows with the new variables.
$ws.Cells.Item(237, 4).Value = "world.real_weather"
$ws.Cells.Item(237, 5).Value = "string"
$ws.Cells.Item(237, 6).Value = "Current weather's name."

$ws.Cells.Item(238, 4).Value = "world.real_weather_element"
$ws.Cells.Item(238, 5).Value = "string"
$ws.Cells.Item(238, 6).Value = "Current weather's element."

# The "World" section header was un-merged in the edit; split it back into
# individually-styled cells (left/general aligned instead of centered).
$ws.Range("A230:A243").UnMerge()
$ws.Range("A230:A243").HorizontalAlignment = 1

# Restore the on-screen selection to roughly where the edit was made.
$ws.Range("E239").Select()
